# Fix "cantidad" column (C) values that had an extra trailing zero
# erroneously appended. Divide each affected value by 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the last used row based on column A (nemotecnico)
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -ne $null) {
        $cell.Value2 = $cell.Value2 / 10
    }
}
